$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of case data appended below the existing table (row 17).
# Columns B (and similar date-looking text columns) must stay as literal
# text rather than being auto-parsed into Excel date serials, so we format
# the cell as Text before assigning the value. Afterwards we reset the
# cell style back to "Normal" so no extra style index is introduced and the
# cell matches the plain (unstyled) look of the other data rows above it.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("A17") "YXYZH"
Set-TextValue $ws.Range("B17") "11/29/2024"
Set-TextValue $ws.Range("C17") "Ashar Nadeem"
Set-TextValue $ws.Range("D17") "0322-7287568"
Set-TextValue $ws.Range("E17") "Lahore"
Set-TextValue $ws.Range("F17") "0322-7287568"
Set-TextValue $ws.Range("G17") "xxdxsdxdsx"
Set-TextValue $ws.Range("H17") "Iqbal town"
Set-TextValue $ws.Range("I17") "knj"
